$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns F, D, C, A (right-to-left so earlier deletions don't
# shift the addresses of columns we still need to delete).
$ws.Range("F1").EntireColumn.Delete()
$ws.Range("D1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()
$ws.Range("A1").EntireColumn.Delete()

# Update the active selection to match the saved view state.
$ws.Range("E7").Select()
